# Updates cryptos list figures (Price and Volume(1h) columns) to match
# the latest scrape, per commit "Updated cryptos list on Tue Jan 23
# 06:40:04 UTC 2024 with GitHub Actions".
#
# The Price/Volume(1h) cells are plain text (inline strings) in the
# workbook, but many of the new values look like ordinary decimal
# numbers (e.g. "2.20", "93.20", "0.0808"). Assigning such strings
# directly via Range.Value causes Excel to auto-convert them to numeric
# values, which silently drops significant trailing zeros and changes
# the stored representation/cell style. To avoid that, each value is
# written as a literal text formula (="...") and then converted in place
# to a plain value via Copy + PasteSpecial (values only); this preserves
# the exact text (including padding spaces on the percentage cells) and
# keeps the original "General" cell formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $escaped = $val.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy() | Out-Null
    $c.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

$excel.CutCopyMode = 0

Set-TextValue 'D2' '40.126.04'
Set-TextValue 'E2' '  -2.53%  '
Set-TextValue 'D3' '2.348.19'
Set-TextValue 'E3' '  -3.37%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '311.32'
Set-TextValue 'E5' '  -2.04%  '
Set-TextValue 'D6' '85.24'
Set-TextValue 'E6' '  -5.04%  '
Set-TextValue 'D7' '0.527'
Set-TextValue 'E7' '  -2.25%  '
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'E9' '  -3.45%  '
Set-TextValue 'D10' '0.0808'
Set-TextValue 'E10' '  -3.62%  '
Set-TextValue 'D11' '30.15'
Set-TextValue 'E11' '  -6.35%  '
Set-TextValue 'E12' '  +0.73%  '
Set-TextValue 'D13' '2.710.21'
Set-TextValue 'E13' '  -3.30%  '
Set-TextValue 'D14' '6.42'
Set-TextValue 'E14' '  -4.63%  '
Set-TextValue 'D15' '14.84'
Set-TextValue 'E15' '  -5.55%  '
Set-TextValue 'D16' '2.365.66'
Set-TextValue 'E16' '  -2.65%  '
Set-TextValue 'D17' '0.762'
Set-TextValue 'E17' '  -1.76%  '
Set-TextValue 'D18' '40.131.44'
Set-TextValue 'E18' '  -2.32%  '
Set-TextValue 'E19' '  -2.85%  '
Set-TextValue 'E20' '  -2.95%  '
Set-TextValue 'D21' '68.30'
Set-TextValue 'E21' '  -4.61%  '
Set-TextValue 'D22' '10.64'
Set-TextValue 'E22' '  -4.34%  '
Set-TextValue 'D23' '235.45'
Set-TextValue 'E23' '  -0.24%  '
Set-TextValue 'E24' '  -5.24%  '
Set-TextValue 'E25' '  -0.02%  '
Set-TextValue 'D26' '1.82'
Set-TextValue 'E26' '  -2.96%  '
Set-TextValue 'D27' '23.67'
Set-TextValue 'E27' '  -2.18%  '
Set-TextValue 'D28' '2.20'
Set-TextValue 'E28' '  -1.54%  '
Set-TextValue 'D29' '9.27'
Set-TextValue 'E29' '  -3.63%  '
Set-TextValue 'D30' '34.86'
Set-TextValue 'E30' '  +0.38%  '
Set-TextValue 'D31' '153.99'
Set-TextValue 'E31' '  -1.04%  '
Set-TextValue 'E32' '  -0.05%  '
Set-TextValue 'D33' '5.11'
Set-TextValue 'E33' '  -3.27%  '
Set-TextValue 'D34' '2.49'
Set-TextValue 'E34' '  -0.83%  '
Set-TextValue 'E35' '  -3.91%  '
Set-TextValue 'E36' '  -0.73%  '
Set-TextValue 'E37' '  -5.32%  '
Set-TextValue 'D38' '0.0993'
Set-TextValue 'E38' '  -1.86%  '
Set-TextValue 'D39' '15.71'
Set-TextValue 'E39' '  -6.46%  '
Set-TextValue 'D40' '1.72'
Set-TextValue 'E40' '  -3.56%  '
Set-TextValue 'D41' '3.85'
Set-TextValue 'E41' '  -1.86%  '
Set-TextValue 'D42' '1.969.07'
Set-TextValue 'E42' '  -1.49%  '
Set-TextValue 'E43' '  -0.08%  '
Set-TextValue 'D44' '0.0266'
Set-TextValue 'E44' '  -3.77%  '
Set-TextValue 'D45' '17.66'
Set-TextValue 'E45' '  -5.20%  '
Set-TextValue 'E46' '  +0.66%  '
Set-TextValue 'D47' '2.69'
Set-TextValue 'E47' '  -7.37%  '
Set-TextValue 'D48' '2.569.54'
Set-TextValue 'E48' '  -3.47%  '
Set-TextValue 'D49' '93.20'
Set-TextValue 'E49' '  -1.89%  '
Set-TextValue 'D50' '70.37'
Set-TextValue 'E50' '  -4.36%  '
Set-TextValue 'D51' '50.25'
Set-TextValue 'E51' '  -3.58%  '

$excel.CutCopyMode = 0
